# The deck ships with two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colour scheme (the one
#                            actually applied to the slide master, i.e. the
#                            design used by every slide in the deck)
# The edit swaps the two designs: the deck-wide design goes from the
# "Integral" (Red Violet) look back to the default "Office Theme" colours.
# We reproduce that by rewriting the live theme's 12 colour-scheme slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the stock Office
# palette, via the ThemeColorScheme exposed on a slide that uses that
# master/theme. Because every slide shares the one slide master, touching
# it once re-colours the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (slot, target Office sRGB) ; RGB property takes a BGR-packed
# long (0x00BBGGRR), same convention PowerPoint COM always uses.
$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
